# Adds columns I (I0) and J (IF) to the worksheet, mirroring column H (IP) pattern
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: I1 = "I0", J1 = "IF"; copy style from existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows 2-66 for columns I and J
$data = @(
    @(2,9,9),
    @(3,7,7),
    @(4,8,8),
    @(5,9,9),
    @(6,6,6),
    @(7,6,6),
    @(8,8,8),
    @(9,8,8),
    @(10,9,9),
    @(11,9,9),
    @(12,9,9),
    @(13,9,9),
    @(14,9,9),
    @(15,9,9),
    @(16,8,8),
    @(17,7,7),
    @(18,9,9),
    @(19,7,7),
    @(20,9,9),
    @(21,9,9),
    @(22,10,10),
    @(23,9,9),
    @(24,9,9),
    @(25,9,9),
    @(26,9,9),
    @(27,9,9),
    @(28,9,9),
    @(29,9,9),
    @(30,9,9),
    @(31,9,9),
    @(32,9,9),
    @(33,9,9),
    @(34,9,9),
    @(35,9,9),
    @(36,9,9),
    @(37,9,9),
    @(38,9,9),
    @(39,9,9),
    @(40,9,9),
    @(41,9,9),
    @(42,9,9),
    @(43,8,9),
    @(44,9,9),
    @(45,9,9),
    @(46,9,9),
    @(47,8,9),
    @(48,8,8),
    @(49,9,9),
    @(50,9,9),
    @(51,9,9),
    @(52,9,9),
    @(53,9,9),
    @(54,9,9),
    @(55,9,9),
    @(56,7,7),
    @(57,9,9),
    @(58,9,9),
    @(59,9,9),
    @(60,9,10),
    @(61,8,8),
    @(62,7,7),
    @(63,9,9),
    @(64,7,7),
    @(65,5,6),
    @(66,4,4)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
